$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 361-380: columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg) ---

$ws.Cells.Item(361, 4).Value = 44753
$ws.Cells.Item(361, 10).Value = 130
$ws.Cells.Item(361, 11).Value = 8000
$ws.Cells.Item(361, 12).Value = 9000
$ws.Cells.Item(361, 13).Value = 8500
$ws.Cells.Item(361, 16).Value = 121

$ws.Cells.Item(362, 4).Value = 44753
$ws.Cells.Item(362, 10).Value = 120
$ws.Cells.Item(362, 11).Value = 6000
$ws.Cells.Item(362, 12).Value = 7000
$ws.Cells.Item(362, 13).Value = 6500
$ws.Cells.Item(362, 16).Value = 65

$ws.Cells.Item(363, 4).Value = 44489
$ws.Cells.Item(363, 10).Value = 120
$ws.Cells.Item(363, 11).Value = 7000
$ws.Cells.Item(363, 12).Value = 8000
$ws.Cells.Item(363, 13).Value = 7500
$ws.Cells.Item(363, 16).Value = 107

$ws.Cells.Item(364, 4).Value = 44489
$ws.Cells.Item(364, 10).Value = 120
$ws.Cells.Item(364, 11).Value = 6000
$ws.Cells.Item(364, 12).Value = 7000
$ws.Cells.Item(364, 13).Value = 6500
$ws.Cells.Item(364, 16).Value = 65

$ws.Cells.Item(365, 4).Value = 44659
$ws.Cells.Item(365, 10).Value = 160
$ws.Cells.Item(365, 11).Value = 4000
$ws.Cells.Item(365, 12).Value = 5000
$ws.Cells.Item(365, 13).Value = 4500
$ws.Cells.Item(365, 16).Value = 64

$ws.Cells.Item(366, 4).Value = 44659
$ws.Cells.Item(366, 10).Value = 170
$ws.Cells.Item(366, 11).Value = 3500
$ws.Cells.Item(366, 12).Value = 4000
$ws.Cells.Item(366, 13).Value = 3750
$ws.Cells.Item(366, 16).Value = 38

$ws.Cells.Item(367, 4).Value = 44340
$ws.Cells.Item(367, 10).Value = 120
$ws.Cells.Item(367, 11).Value = 9000
$ws.Cells.Item(367, 12).Value = 10000
$ws.Cells.Item(367, 13).Value = 9500
$ws.Cells.Item(367, 16).Value = 136

$ws.Cells.Item(368, 4).Value = 44340
$ws.Cells.Item(368, 10).Value = 120
$ws.Cells.Item(368, 11).Value = 8000
$ws.Cells.Item(368, 12).Value = 9000
$ws.Cells.Item(368, 13).Value = 8500
$ws.Cells.Item(368, 16).Value = 85

$ws.Cells.Item(369, 4).Value = 44326
$ws.Cells.Item(369, 10).Value = 120
$ws.Cells.Item(369, 11).Value = 9000
$ws.Cells.Item(369, 12).Value = 10000
$ws.Cells.Item(369, 13).Value = 9500
$ws.Cells.Item(369, 16).Value = 136

$ws.Cells.Item(370, 4).Value = 44326
$ws.Cells.Item(370, 10).Value = 120
$ws.Cells.Item(370, 11).Value = 8000
$ws.Cells.Item(370, 12).Value = 9000
$ws.Cells.Item(370, 13).Value = 8500
$ws.Cells.Item(370, 16).Value = 85

$ws.Cells.Item(371, 4).Value = 44343
$ws.Cells.Item(371, 10).Value = 120
$ws.Cells.Item(371, 11).Value = 8500
$ws.Cells.Item(371, 12).Value = 9000
$ws.Cells.Item(371, 13).Value = 8750
$ws.Cells.Item(371, 16).Value = 125

$ws.Cells.Item(372, 4).Value = 44343
$ws.Cells.Item(372, 10).Value = 120
$ws.Cells.Item(372, 11).Value = 7500
$ws.Cells.Item(372, 12).Value = 8000
$ws.Cells.Item(372, 13).Value = 7750
$ws.Cells.Item(372, 16).Value = 78

$ws.Cells.Item(373, 4).Value = 44426
$ws.Cells.Item(373, 10).Value = 120
$ws.Cells.Item(373, 11).Value = 11000
$ws.Cells.Item(373, 12).Value = 12000
$ws.Cells.Item(373, 13).Value = 11500
$ws.Cells.Item(373, 16).Value = 164

$ws.Cells.Item(374, 4).Value = 44426
$ws.Cells.Item(374, 10).Value = 120
$ws.Cells.Item(374, 11).Value = 8000
$ws.Cells.Item(374, 12).Value = 9000
$ws.Cells.Item(374, 13).Value = 8500
$ws.Cells.Item(374, 16).Value = 85

$ws.Cells.Item(375, 4).Value = 44259
$ws.Cells.Item(375, 10).Value = 120
$ws.Cells.Item(375, 11).Value = 8000
$ws.Cells.Item(375, 12).Value = 9000
$ws.Cells.Item(375, 13).Value = 8500
$ws.Cells.Item(375, 16).Value = 121

$ws.Cells.Item(376, 4).Value = 44259
$ws.Cells.Item(376, 10).Value = 120
$ws.Cells.Item(376, 11).Value = 6000
$ws.Cells.Item(376, 12).Value = 7000
$ws.Cells.Item(376, 13).Value = 6500
$ws.Cells.Item(376, 16).Value = 65

$ws.Cells.Item(377, 4).Value = 44376
$ws.Cells.Item(377, 10).Value = 190
$ws.Cells.Item(377, 11).Value = 8000
$ws.Cells.Item(377, 12).Value = 9000
$ws.Cells.Item(377, 13).Value = 8684
$ws.Cells.Item(377, 16).Value = 124

$ws.Cells.Item(378, 4).Value = 44376
$ws.Cells.Item(378, 10).Value = 120
$ws.Cells.Item(378, 11).Value = 6000
$ws.Cells.Item(378, 12).Value = 7000
$ws.Cells.Item(378, 13).Value = 6500
$ws.Cells.Item(378, 16).Value = 65

$ws.Cells.Item(379, 4).Value = 44627
$ws.Cells.Item(379, 10).Value = 160
$ws.Cells.Item(379, 11).Value = 7000
$ws.Cells.Item(379, 12).Value = 7500
$ws.Cells.Item(379, 13).Value = 7250
$ws.Cells.Item(379, 16).Value = 104

$ws.Cells.Item(380, 4).Value = 44627
$ws.Cells.Item(380, 10).Value = 120
$ws.Cells.Item(380, 11).Value = 6000
$ws.Cells.Item(380, 12).Value = 6500
$ws.Cells.Item(380, 13).Value = 6250
$ws.Cells.Item(380, 16).Value = 62

# --- Append new rows 381 and 382 (same master data, shifted-in week) ---

$ws.Cells.Item(381, 1).Value = 1
$ws.Cells.Item(381, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(381, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(381, 4).Value = 44454
$ws.Cells.Item(381, 5).Value = 15
$ws.Cells.Item(381, 6).Value = 100112032
$ws.Cells.Item(381, 7).Value = "Zapallo italiano"
$ws.Cells.Item(381, 8).Value = "Huracán"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 130
$ws.Cells.Item(381, 11).Value = 6000
$ws.Cells.Item(381, 12).Value = 7000
$ws.Cells.Item(381, 13).Value = 6500
$ws.Cells.Item(381, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(381, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(381, 16).Value = 93
$ws.Cells.Item(381, 17).Value = 70
$ws.Cells.Item(381, 18).Value = "Hortaliza"
$ws.Cells.Item(381, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(382, 1).Value = 1
$ws.Cells.Item(382, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(382, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(382, 4).Value = 44454
$ws.Cells.Item(382, 5).Value = 15
$ws.Cells.Item(382, 6).Value = 100112032
$ws.Cells.Item(382, 7).Value = "Zapallo italiano"
$ws.Cells.Item(382, 8).Value = "Huracán"
$ws.Cells.Item(382, 9).Value = "Segunda"
$ws.Cells.Item(382, 10).Value = 136
$ws.Cells.Item(382, 11).Value = 4000
$ws.Cells.Item(382, 12).Value = 5000
$ws.Cells.Item(382, 13).Value = 4500
$ws.Cells.Item(382, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 45
$ws.Cells.Item(382, 17).Value = 100
$ws.Cells.Item(382, 18).Value = "Hortaliza"
$ws.Cells.Item(382, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

